$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.174.36"
$ws.Range("E2").Value = "  -1.83%  "
$ws.Range("D3").Value = "1.659.64"
$ws.Range("E3").Value = "  -1.81%  "
$ws.Range("E4").Value = "  +0.49%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.65"
$ws.Range("E5").Value = "  -0.07%  "
$ws.Range("E6").Value = "  -2.34%  "
$ws.Range("E7").Value = "  +0.50%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2664"
$ws.Range("E8").Value = "  -0.78%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06317"
$ws.Range("E9").Value = "  -1.82%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.04"
$ws.Range("E10").Value = "  -2.94%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07724"
$ws.Range("E11").Value = "  -0.94%  "
$ws.Range("D12").Value = "1.678.93"
$ws.Range("E12").Value = "  -0.54%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.426"
$ws.Range("E13").Value = "  -1.77%  "
$ws.Range("D14").Value = "1.888.18"
$ws.Range("E14").Value = "  -1.68%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5450"
$ws.Range("E15").Value = "  -3.17%  "
$ws.Range("D16").Value = "0.0₅8214"
$ws.Range("E16").Value = "  -2.83%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.83"
$ws.Range("E17").Value = "  -2.35%  "
$ws.Range("D18").Value = "26.216.47"
$ws.Range("E18").Value = "  -1.82%  "
$ws.Range("E19").Value = "  +0.41%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.656"
$ws.Range("E20").Value = "  -3.36%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "192.94"
$ws.Range("E21").Value = "  -1.49%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.14"
$ws.Range("E22").Value = "  -2.77%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.078"
$ws.Range("E23").Value = "  -4.84%  "
$ws.Range("E24").Value = "  +0.65%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "138.82"
$ws.Range("E25").Value = "  -3.68%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1236"
$ws.Range("E26").Value = "  -4.50%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.205"
$ws.Range("E27").Value = "  -3.84%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.13"
$ws.Range("E28").Value = "  -0.92%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.415"
$ws.Range("E29").Value = "  -0.81%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.06002"
$ws.Range("E30").Value = "  -2.79%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.282"
$ws.Range("E31").Value = "  +0.00%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.596"
$ws.Range("E32").Value = "  -0.24%  "
$ws.Range("E33").Value = "  -4.15%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.643"
$ws.Range("E34").Value = "  -3.60%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9802"
$ws.Range("E35").Value = "  -3.46%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.786"
$ws.Range("E36").Value = "  -0.53%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.413"
$ws.Range("E37").Value = "  -0.40%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5910"
$ws.Range("E38").Value = "  +2.91%  "
$ws.Range("E39").Value = "  -3.95%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.949"
$ws.Range("E40").Value = "  -1.23%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8656"
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.005"
$ws.Range("E42").Value = "  +0.38%  "
$ws.Range("D43").Value = "1.037.64"
$ws.Range("E43").Value = "  -3.98%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.63"
$ws.Range("E44").Value = "  -0.87%  "
$ws.Range("D45").Value = "1.802.88"
$ws.Range("E45").Value = "  -2.04%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "57.06"
$ws.Range("E46").Value = "  -0.80%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "0.0₈107"
$ws.Range("E47").Value = "  -2.29%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.004"
$ws.Range("E48").Value = "  +0.20%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.109"
$ws.Range("E49").Value = "  -1.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05186"
$ws.Range("E50").Value = "  -0.72%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4232"
$ws.Range("E51").Value = "  -0.20%  "
